# Apply updated cryptos list values (price + 1h volume change) to Sheet1.
# Cells whose new text is a bare number get a leading apostrophe so Excel
# keeps storing them as text (matching the source data's inline-string cells)
# instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.338.41"
$ws.Range("E2").Value = "  +5.42%  "
$ws.Range("D3").Value = "4.076.97"
$ws.Range("E3").Value = "  +6.00%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'522.35"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'148.70"
$ws.Range("E6").Value = "  +3.84%  "
$ws.Range("D7").Value = "'0.724"
$ws.Range("E7").Value = "  +19.78%  "
$ws.Range("D8").Value = "4.070.31"
$ws.Range("E8").Value = "  +5.96%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "'0.781"
$ws.Range("E10").Value = "  +9.97%  "
$ws.Range("D11").Value = "'0.179"
$ws.Range("E11").Value = "  +6.09%  "
$ws.Range("D12").Value = "'0.0000333"
$ws.Range("E12").Value = "  +2.23%  "
$ws.Range("D13").Value = "'48.64"
$ws.Range("E13").Value = "  +16.92%  "
$ws.Range("D14").Value = "'11.06"
$ws.Range("E14").Value = "  +9.07%  "
$ws.Range("D15").Value = "4.715.30"
$ws.Range("E15").Value = "  +5.94%  "
$ws.Range("D16").Value = "4.070.18"
$ws.Range("E16").Value = "  +5.62%  "
$ws.Range("D17").Value = "'14.57"
$ws.Range("E17").Value = "  +5.10%  "
$ws.Range("D18").Value = "'21.37"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "'1.25"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "72.442.38"
$ws.Range("E21").Value = "  +5.47%  "
$ws.Range("D22").Value = "'448.36"
$ws.Range("E22").Value = "  +7.06%  "
$ws.Range("D23").Value = "'103.76"
$ws.Range("E23").Value = "  +20.01%  "
$ws.Range("D24").Value = "'3.62"
$ws.Range("D25").Value = "'15.08"
$ws.Range("E25").Value = "  +7.90%  "
$ws.Range("D26").Value = "'4.05"
$ws.Range("E26").Value = "  +2.49%  "
$ws.Range("D27").Value = "'11.47"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").Value = "'11.15"
$ws.Range("E28").Value = "  +5.92%  "
$ws.Range("D29").Value = "'38.14"
$ws.Range("E29").Value = "  +6.18%  "
$ws.Range("E30").Value = "  +2.99%  "
$ws.Range("D31").Value = "'3.30"
$ws.Range("E31").Value = "  +16.73%  "
$ws.Range("D32").Value = "'13.72"
$ws.Range("E32").Value = "  +5.04%  "
$ws.Range("D33").Value = "'0.132"
$ws.Range("E33").Value = "  +5.47%  "
$ws.Range("D34").Value = "'684.02"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("D35").Value = "'67.77"
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").Value = "'6.67"
$ws.Range("E36").Value = "  +14.05%  "
$ws.Range("D37").Value = "'42.38"
$ws.Range("E37").Value = "  +7.09%  "
$ws.Range("D38").Value = "0.0₃0876"
$ws.Range("E38").Value = "  +3.48%  "
$ws.Range("D39").Value = "'0.433"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").Value = "'0.154"
$ws.Range("E40").Value = "  +5.45%  "
$ws.Range("D41").Value = "'3.47"
$ws.Range("E41").Value = "  +9.92%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").Value = "'0.0504"
$ws.Range("E43").Value = "  +6.00%  "
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "'3.23"
$ws.Range("E45").Value = "  +2.78%  "
$ws.Range("D46").Value = "'0.158"
$ws.Range("E46").Value = "  +14.13%  "
$ws.Range("D47").Value = "'9.91"
$ws.Range("E47").Value = "  +17.84%  "
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").Value = "'3.43"
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D50").Value = "'3.08"
$ws.Range("E50").Value = "  +5.18%  "
$ws.Range("D51").Value = "'0.000280"
$ws.Range("E51").Value = "  +3.69%  "
